$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure rows 2:3 exist as real (present-but-empty) cells across A:F,
# matching the "templated" row layout, by touching their style first.
$ws.Range("A2:F3").Style = "Normal"

# Populate the "creation_method" column (C) with the two possible values.
$ws.Range("C2").Value = "translational"
$ws.Range("C3").Value = "rotational"
